# Saved progress at the end of the loop
# Re-derive the billing detail rows: drop the old "Short point" line item
# (row 9) which shifts every row below it up by one, then refresh the
# quantities / computed amounts for the remaining line items and the
# Grand Total / Net Payable summary rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Short point (up to 3 mtr.)" row entirely; everything below
# (rows 10-20) shifts up to (rows 9-19).
$ws.Range("A9").EntireRow.Delete()

# --- Row 8: "Rewiring of light point/ fan point/..." quantity refreshed ---
$ws.Range("C8").Value = 2

# --- Row 9 (was row 10 "Long point"): now "Medium point (up to 6 mtr.)" ---
$ws.Range("C9").Value = 85
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "Medium point (up to 6 mtr.)"
$ws.Range("F9").Value = 472
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "40120.00"
$ws.Range("G9").Style = "Normal"

# --- Row 10 (was row 11 "Rewiring of 3/5 pin..."): quantity refreshed ---
$ws.Range("C10").Value = 91

# --- Row 11 (was row 12 "On board"): quantity + upto-date amount refreshed ---
$ws.Range("C11").Value = 28
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "3808.00"
$ws.Range("G11").Style = "Normal"

# --- Row 12 (was row 13 "P & F ISI marked..."): quantity + upto-date amount refreshed ---
$ws.Range("C12").Value = 32
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "736.00"
$ws.Range("G12").Style = "Normal"

# --- Row 13 (was row 14 "Total"): quantity refreshed ---
$ws.Range("C13").Value = 36

# --- Row 14 (was row 15 "Add Tender Premium"): quantity refreshed ---
$ws.Range("C14").Value = 5

# --- Row 15 (was row 16 "Grand Total"): quantity refreshed ---
$ws.Range("C15").Value = 26

# --- Row 17 "Grand Total Rs." summary amounts refreshed ---
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "44664.00"
$ws.Range("G17").Style = "Normal"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "44664.00"
$ws.Range("H17").Style = "Normal"

# --- Row 19 "NET PAYABLE AMOUNT Rs." summary amounts refreshed ---
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "44664.00"
$ws.Range("G19").Style = "Normal"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "44664.00"
$ws.Range("H19").Style = "Normal"
